$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-404.6565856933594, -0.3036, -0.1298000067472458, 0.4621999859809875, -1.169700026512146)
    3 = @(-165.6318054199219, -0.1325, -0.0462, 0.7084000110626221, -1.012199997901917)
    4 = @(-64.75749969482422, -0.0757, -0.0468, 0.6057999730110168, -0.5006999969482422)
    5 = @(111.1718978881836, 0.1173, 0.019, 1.068300008773804, -0.3077999949455261)
    6 = @(45.53419876098633, 0.036, -0.0155, 0.9750999808311462, -0.4165999889373779)
    7 = @(-98.08029937744141, -0.0988, -0.0706000030040741, 0.5040000081062317, -0.4607000052928925)
    8 = @(-160.0614013671875, -0.1434, -0.1189, 0.2356999963521957, -0.6043999791145325)
    9 = @(-736.4815063476562, -0.0954, -0.0452, 1.068300008773804, -1.169700026512146)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
}
